$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2994946666666667
$ws.Range("H2").Value = 0.8984840000000001
$ws.Range("J2").Value = 0.4989451716962828
$ws.Range("M2").Value = 0.07271233333333334
$ws.Range("N2").Value = 0.218137
$ws.Range("O2").Value = 0.004171225362010892
$ws.Range("P2").Value = 0.004171225362010893
$ws.Range("Q2").Value = 0.02177695603422223
$ws.Range("R2").Value = 0.195992604308
$ws.Range("S2").Value = 0.002081212754432414
$ws.Range("T2").Value = 0.002081212754432414

$ws.Range("G3").Value = 0.2994946666666667
$ws.Range("H3").Value = 0.8984840000000001
$ws.Range("J3").Value = 0.4989451716962828
$ws.Range("O3").Value = 0.5387060579248023
$ws.Range("P3").Value = 0.5387060579248023
$ws.Range("Q3").Value = 2.812453684627112
$ws.Range("S3").Value = 0.2687847865651181
$ws.Range("T3").Value = 0.2687847865651181

$ws.Range("G4").Value = 0.2994946666666667
$ws.Range("H4").Value = 0.8984840000000001
$ws.Range("J4").Value = 0.4989451716962828
$ws.Range("O4").Value = 0.4571227167131868
$ws.Range("P4").Value = 0.4571227167131868
$ws.Range("Q4").Value = 2.386526845269334
$ws.Range("S4").Value = 0.2280791723767322
$ws.Range("T4").Value = 0.2280791723767322

$ws.Range("I5").Value = 0.5010548283037172
$ws.Range("M5").Value = 0.07271233333333334
$ws.Range("N5").Value = 0.218137
$ws.Range("O5").Value = 0.004171225362010892
$ws.Range("P5").Value = 0.004171225362010893
$ws.Range("Q5").Value = 0.02186903408566667
$ws.Range("R5").Value = 0.196821306771
$ws.Range("S5").Value = 0.002090012607578478
$ws.Range("T5").Value = 0.002090012607578478

$ws.Range("I6").Value = 0.5010548283037172
$ws.Range("O6").Value = 0.5387060579248023
$ws.Range("P6").Value = 0.5387060579248023
$ws.Range("S6").Value = 0.2699212713596841
$ws.Range("T6").Value = 0.2699212713596841

$ws.Range("I7").Value = 0.5010548283037172
$ws.Range("O7").Value = 0.4571227167131868
$ws.Range("P7").Value = 0.4571227167131868
$ws.Range("S7").Value = 0.2290435443364545
$ws.Range("T7").Value = 0.2290435443364545
